$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.342.32'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").Value = '1.841.73'
$ws.Range("E3").Value = '  -0.20%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9989'
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.54'
$ws.Range("E5").Value = '  -0.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6283'
$ws.Range("E6").Value = '  -0.47%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07438'
$ws.Range("E8").Value = '  -0.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2891'
$ws.Range("E9").Value = '  -0.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.93'
$ws.Range("E10").Value = '  +2.13%  '

$ws.Range("E11").Value = '  +0.04%  '

$ws.Range("D12").Value = '1.842.51'
$ws.Range("E12").Value = '  -0.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.963'
$ws.Range("E13").Value = '  -0.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6753'
$ws.Range("E14").Value = '  -0.49%  '

$ws.Range("E15").Value = '  +0.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.48'
$ws.Range("E16").Value = '  -0.81%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.233'
$ws.Range("E17").Value = '  +1.47%  '

$ws.Range("D18").Value = '29.389.59'
$ws.Range("E18").Value = '  -0.09%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.80'
$ws.Range("E19").Value = '  +0.01%  '

$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("E22").Value = '  -1.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.04'
$ws.Range("E24").Value = '  -0.56%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.489'
$ws.Range("E25").Value = '  +0.76%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1346'
$ws.Range("E26").Value = '  -2.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.40'
$ws.Range("E27").Value = '  -0.82%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.07086'
$ws.Range("E28").Value = '  +11.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.463'
$ws.Range("E29").Value = '  +5.83%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.478'
$ws.Range("E30").Value = '  +0.26%  '

$ws.Range("E31").Value = '  -1.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.032'
$ws.Range("E32").Value = '  -0.57%  '

$ws.Range("E33").Value = '  +0.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.138'
$ws.Range("E34").Value = '  -0.28%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6968'
$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.581'

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01839'
$ws.Range("E37").Value = '  +0.93%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.811'
$ws.Range("E38").Value = '  -0.87%  '

$ws.Range("D39").Value = '1.235.41'
$ws.Range("E39").Value = '  -1.46%  '

$ws.Range("E40").Value = '  +3.79%  '

$ws.Range("E41").Value = '  +2.19%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  +0.10%  '

$ws.Range("D43").Value = '2.002.00'
$ws.Range("E43").Value = '  -0.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.76'
$ws.Range("E44").Value = '  -0.61%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.25'
$ws.Range("E45").Value = '  -1.65%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000120'
$ws.Range("E46").Value = '  +1.02%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.019'
$ws.Range("E47").Value = '  -0.42%  '

$ws.Range("E48").Value = '  +1.07%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.927'
$ws.Range("E49").Value = '  -1.28%  '

$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1139'
$ws.Range("E50").Value = '  -2.94%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3910'
$ws.Range("E51").Value = '  -0.74%  '
